# Update the "Gaz" sheet (2nd worksheet) of the EPEX spot prices workbook:
# prepend a new day's row (2025-06-16, all "-") above the existing
# 2025-06-17 row, pushing it down to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a fresh row above the current row 2 (the 2025-06-17 data row),
# shifting it down to row 3 and growing the sheet's used range to A1:D3.
$ws.Rows.Item(2).Insert()

# Write the new day's values. The leading apostrophe forces text
# interpretation so "2025-06-16" is stored as a literal string instead of
# being auto-converted to a date serial by Excel's smart typing.
$ws.Range("A2").Value = "'2025-06-16"
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"

# The inserted row inherits formatting from the row above (the bold header
# row); strip it back to the default/unstyled look used by the rest of the
# data rows (matches row 3, e.g. the former 2025-06-17 row).
$ws.Range("A2:D2").ClearFormats()
